$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Flag" values for rows 3 through 31 from "No" to "Yes"
for ($r = 3; $r -le 31; $r++) {
    $ws.Range("E$r").Value = "Yes"
}

# Update the current selection to match the new active range (E2:E31, active cell E2)
$ws.Activate()
$ws.Range("E2:E31").Select()
$excel.ActiveCell
